# Generate Report for Handoff
# Refresh the "Latest Handoff Date" / "Latest Handoff Datetime" values for the
# file 7218e061-032c-4c7d-89c4-fa7b3938d44c (row 6 in each sheet) to reflect a
# newly generated handoff.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview!D6 - "Latest Handoff Date" column for 7218e061-032c-4c7d-89c4-fa7b3938d44c.md
$overview.Range("D6").Value = "2016-03-22 14:42:09"

# zh-cn!E6 - "Latest Handoff Datetime" column for the same file
$zhcn.Range("E6").Value = "2016-03-22 14:42:04"

# de-de!E6 - "Latest Handoff Datetime" column for the same file
$dede.Range("E6").Value = "2016-03-22 14:42:09"
